$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 945
$ws1.Range("F3").Value = 1038
$ws1.Range("F5").Value = 896
$ws1.Range("F6").Value = 474
$ws1.Range("F9").Value = 1339
$ws1.Range("F10").Value = 774
$ws1.Range("F12").Value = 577
$ws1.Range("F13").Value = 194
$ws1.Range("F14").Value = 79
$ws1.Range("F15").Value = 79
$ws1.Range("F16").Value = 1345
$ws1.Range("F17").Value = 156
$ws1.Range("F18").Value = 22
$ws1.Range("F20").Value = 18
$ws1.Range("F24").Value = 174
$ws1.Range("F25").Value = 677
$ws1.Range("F27").Value = 1195
$ws1.Range("F28").Value = 32
$ws1.Range("F29").Value = 19

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 648
$ws2.Range("F7").Value = 263

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 945
$ws4.Range("F5").Value = 1038
$ws4.Range("F7").Value = 896
$ws4.Range("F8").Value = 474
$ws4.Range("F11").Value = 1339
$ws4.Range("F12").Value = 774
$ws4.Range("F16").Value = 577
$ws4.Range("F17").Value = 648
$ws4.Range("F18").Value = 194
$ws4.Range("F19").Value = 79
$ws4.Range("F20").Value = 79
$ws4.Range("F21").Value = 1345
$ws4.Range("F23").Value = 156
$ws4.Range("F24").Value = 22
$ws4.Range("F26").Value = 18
$ws4.Range("F29").Value = 263
$ws4.Range("F36").Value = 174
$ws4.Range("F37").Value = 677
$ws4.Range("F39").Value = 1195
$ws4.Range("F40").Value = 32
$ws4.Range("F41").Value = 19
